# nexial-macro.xlsx: add new "web" commands (check/checkByLocator/uncheck/uncheckByLocator into
# the "desktop" helper list, openInTab into the "web" helper list) on the hidden '#system' sheet.
# These helper columns back named ranges (desktop, web, ...) consumed by the MacroLibrary sheet's
# data-validation dropdowns, so inserting into the middle of a column must shift only that column
# (not the whole row) and the corresponding named range must be resized to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---- helper: shift a single-column block of cells down by $n rows (keeps other columns intact) ----
function Shift-ColumnDown {
    param($sheet, $col, $firstRow, $lastRow, $n)
    for ($r = $lastRow; $r -ge $firstRow; $r--) {
        $v = $sheet.Cells.Item($r, $col).Value()
        $sheet.Cells.Item($r + $n, $col).Value = $v
    }
}

# column H = 8 ("desktop" named range), column AA = 27 ("web" named range)
$colH = 8
$colAA = 27

# ---- 1) insert check(name) / checkByLocator(locator) before the old H26 ----
Shift-ColumnDown $ws $colH 26 107 2
$ws.Cells.Item(26, $colH).Value = "check(name)"
$ws.Cells.Item(27, $colH).Value = "checkByLocator(locator)"

# ---- 2) insert uncheck(name) / uncheckByLocator(locator) before what is now H102 ----
Shift-ColumnDown $ws $colH 102 109 2
$ws.Cells.Item(102, $colH).Value = "uncheck(name)"
$ws.Cells.Item(103, $colH).Value = "uncheckByLocator(locator)"

# ---- 3) insert openInTab(name,url) before the old AA85 ----
Shift-ColumnDown $ws $colAA 85 152 1
$ws.Cells.Item(85, $colAA).Value = "openInTab(name,url)"

# ---- 4) resize the named ranges backing the two lists we edited ----
$wb.Names.Item("desktop").RefersTo = "='#system'!`$H`$2:`$H`$111"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$153"
